# Add default gym property
# ----------------------------------------------------------------------
# This script rewrites two paragraphs inside the "ALGORITHM" section of
# the document:
#
#   1. "User tries to check the current gym for it's current ..." -- the
#      run holding "it's" is split out into its own run wrapped in
#      w:proofErr spell-check markers (no textual change).
#
#   2. "IF user is logged in ... DoorCloseTimer ticks down" -- several
#      camel-case identifiers (IsOpenDoorRequested, DoorOpened,
#      DoorCloseTimer) get wrapped in w:proofErr spell-check markers, and
#      new algorithm steps are appended describing what happens when the
#      door-close timer is interrupted / resets.
#
# Because w:proofErr markers and fine-grained run-splits are not
# reproducible through simple Find/Replace, each paragraph's body is
# rebuilt as literal WordprocessingML and dropped in with Range.InsertXML
# -- using the paragraph's Range minus its trailing paragraph-mark
# character keeps the <w:p> element itself (and its paraId/rsid
# attributes) untouched while swapping out only the run content.
# ----------------------------------------------------------------------

$d = $word.ActiveDocument

function Get-ParagraphByText($doc, [string]$needle) {
    $r = $doc.Content
    $r.Find.ClearFormatting()
    $ok = $r.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        return $null
    }
    return $r.Paragraphs(1)
}

function Set-ParagraphInnerXml($doc, $para, [string]$innerXml) {
    $full = $para.Range
    $r = $doc.Range($full.Start, $full.End - 1)
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' `
        + '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' `
        + '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' `
        + '<pkg:xmlData>' `
        + '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' `
        + '<w:body><w:p>' + $innerXml + '</w:p></w:body>' `
        + '</w:document>' `
        + '</pkg:xmlData></pkg:part></pkg:package>'
    $r.InsertXML($pkg)
}

# --- Paragraph 1: "User tries to check the current gym for it's current" ---
$para1 = Get-ParagraphByText $d "User tries to check the current gym for it"

$para1Xml = '<w:r><w:t xml:space="preserve">User tries to check the current gym for </w:t></w:r>' `
    + '<w:proofErr w:type="spellStart"/>' `
    + '<w:r><w:t>it' + [char]0x2019 + 's</w:t></w:r>' `
    + '<w:proofErr w:type="spellEnd"/>' `
    + '<w:r><w:t xml:space="preserve"> current  </w:t></w:r>' `
    + '<w:r><w:t>number of visitors, opens the home page.</w:t></w:r>' `
    + '<w:r><w:br/><w:t xml:space="preserve">IF user is not logged in </w:t></w:r>' `
    + '<w:r><w:br/></w:r>' `
    + '<w:r><w:tab/><w:t>THEN re-direct to login screen.</w:t></w:r>'

Set-ParagraphInnerXml $d $para1 $para1Xml

# --- Paragraph 2: "IF user is logged in ... DoorCloseTimer ticks down" ---
$para2 = Get-ParagraphByText $d "IF user is logged in"

$para2Xml = '<w:r><w:t>IF user is logged in</w:t></w:r>' `
    + '<w:r><w:br/></w:r>' `
    + '<w:r><w:tab/><w:t>THEN send to landing page</w:t></w:r>' `
    + '<w:r><w:br/><w:t>IF User Clicks Enter GYM</w:t></w:r>' `
    + '<w:r><w:br/></w:r>' `
    + '<w:r><w:tab/><w:t xml:space="preserve">THEN </w:t></w:r>' `
    + '<w:proofErr w:type="spellStart"/>' `
    + '<w:r><w:t>IsOpenDoorRequested</w:t></w:r>' `
    + '<w:proofErr w:type="spellEnd"/>' `
    + '<w:r><w:t xml:space="preserve"> is true, </w:t></w:r>' `
    + '<w:proofErr w:type="spellStart"/>' `
    + '<w:r><w:t>DoorOpened</w:t></w:r>' `
    + '<w:proofErr w:type="spellEnd"/>' `
    + '<w:r><w:t xml:space="preserve"> is true</w:t></w:r>' `
    + '<w:r><w:br/><w:t xml:space="preserve"> </w:t></w:r>' `
    + '<w:r><w:tab/></w:r>' `
    + '<w:proofErr w:type="spellStart"/>' `
    + '<w:r><w:t>DoorCloseTimer</w:t></w:r>' `
    + '<w:proofErr w:type="spellEnd"/>' `
    + '<w:r><w:t xml:space="preserve"> starts</w:t></w:r>' `
    + '<w:r><w:br/><w:t xml:space="preserve">              </w:t></w:r>' `
    + '<w:proofErr w:type="spellStart"/>' `
    + '<w:r><w:t>IsOpenDoorRequested</w:t></w:r>' `
    + '<w:proofErr w:type="spellEnd"/>' `
    + '<w:r><w:t xml:space="preserve"> turns false</w:t></w:r>' `
    + '<w:r><w:br/></w:r>' `
    + '<w:r><w:tab/><w:t xml:space="preserve">WHILE </w:t></w:r>' `
    + '<w:proofErr w:type="spellStart"/>' `
    + '<w:r><w:t>DoorCloseTimer</w:t></w:r>' `
    + '<w:proofErr w:type="spellEnd"/>' `
    + '<w:r><w:t xml:space="preserve"> is not 0, </w:t></w:r>' `
    + '<w:r><w:t xml:space="preserve">the </w:t></w:r>' `
    + '<w:proofErr w:type="spellStart"/>' `
    + '<w:r><w:t>DoorCloseTimer</w:t></w:r>' `
    + '<w:proofErr w:type="spellEnd"/>' `
    + '<w:r><w:t xml:space="preserve"> ticks</w:t></w:r>' `
    + '<w:r><w:t xml:space="preserve"> down</w:t></w:r>' `
    + '<w:r><w:br/><w:t xml:space="preserve">              </w:t></w:r>' `
    + '<w:r><w:tab/></w:r>' `
    + '<w:r><w:tab/><w:t xml:space="preserve">IF </w:t></w:r>' `
    + '<w:proofErr w:type="spellStart"/>' `
    + '<w:r><w:t>IsOpenDoorRequested</w:t></w:r>' `
    + '<w:proofErr w:type="spellEnd"/>' `
    + '<w:r><w:t xml:space="preserve"> true while counting down</w:t></w:r>' `
    + '<w:r><w:br/></w:r>' `
    + '<w:r><w:tab/></w:r>' `
    + '<w:r><w:tab/></w:r>' `
    + '<w:proofErr w:type="spellStart"/>' `
    + '<w:r><w:t>DoorCloseTimer</w:t></w:r>' `
    + '<w:proofErr w:type="spellEnd"/>' `
    + '<w:r><w:t xml:space="preserve"> resets</w:t></w:r>' `
    + '<w:r><w:br/></w:r>' `
    + '<w:r><w:tab/></w:r>' `
    + '<w:r><w:tab/><w:t>E</w:t></w:r>'

Set-ParagraphInnerXml $d $para2 $para2Xml

Write-Output "Done."
